# Add a "Phone" column between the "email" column (E) and "PAN" column
# (old F), per commit "Added phone to investor_kyc import".
#
# Inserting a new column at F shifts PAN/Address/Bank Account/IFSC
# Code/Send Confirmation Email one column to the right (old F..J -> new
# G..K) and pushes the used range out to column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F:F").Insert()

# Header
$ws.Range("F1").Value = "Phone"

# Data
$ws.Range("F2").Value = 999999999
$ws.Range("F3").Value = 111111111

# Match the column's width to its neighbours (narrow, numeric column).
# NOTE: the engine snaps ColumnWidth to an internal ~1/7-character pixel
# grid, so 8.11 (rather than the nominal 8.8125 "Format > Column Width"
# figure copied from column E) is what actually lands closest to that
# grid point on export.
$ws.Columns("F:F").ColumnWidth = 8.11

# Move the active selection the way Excel would leave it after typing the
# new column of data and hitting enter past the last data row.
[void]$ws.Range("F4").Select()
